# "veel shit waaronder tier 2 flame tower upgrades"
# Rename the three towers to their "... tower" names, and flesh out the
# Flame tower sheet with its tier-2 -> tier-3 upgrade nodes.

$wb = $excel.ActiveWorkbook

# --- Sheet1: Machine gun -> Machine gun tower -------------------------------
$ws1 = $wb.Worksheets.Item("Machine gun")
$ws1.Range("A2").Value = "Machine gun tower"
$ws1.Range("A2").Select()

# --- Sheet2: Flame thrower -> Flame tower, plus new tier2/tier3 nodes ------
$ws2 = $wb.Worksheets.Item("Flame thrower")
$ws2.Range("A2").Value = "Flame tower"

# existing tier-2 node (D5) grows a tier-3 sibling (F5); two more tier-2
# nodes (D6, D8) are added underneath, each with their own tier-3 sibling.
$ws2.Range("D5").Value = "U0.0.1.4"
$ws2.Range("F5").Value = "U0.0.1.9"
$ws2.Range("D6").Value = "U0.0.1.5"
$ws2.Range("F6").Value = "U0.0.1.10"
$ws2.Range("D8").Value = "U0.0.1.6"
$ws2.Range("F8").Value = "U0.0.1.11"

# the two pre-existing tier-2 nodes (B2 -> D2, B5 -> D3) also gain tier-3
# siblings now that the tier-3 column is in use.
$ws2.Range("F2").Value = "U0.0.1.7"
$ws2.Range("F3").Value = "U0.0.1.8"

# --- Sheet3: Rocket launcher -> Rocket tower --------------------------------
$ws3 = $wb.Worksheets.Item("Rocket launcher")
$ws3.Range("A2").Value = "Rocket tower"

# --- Sheet4: Freezing tower --------------------------------------------------
# (unchanged content; nothing to do)

# --- Rename the sheet tabs themselves ---------------------------------------
$ws1.Name = "Machine gun tower"
$ws2.Name = "Flame tower"
$ws3.Name = "Rocket tower"

# --- Restore the original active sheet / selection on the Flame tower sheet -
$ws2.Range("G6").Select()
